$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.626
$ws.Range("B7").Value = 5.153
$ws.Range("A8").Value = -22.366
$ws.Range("A10").Value = -21.604
$ws.Range("A12").Value = -21.606
$ws.Range("B15").Value = 5.363999999999999
$ws.Range("A18").Value = -21.407
$ws.Range("B18").Value = 6.206
$ws.Range("E18").Value = 16.236
$ws.Range("E19").Value = 16.483
$ws.Range("B20").Value = 6.407000000000001
$ws.Range("E27").Value = 16.392
$ws.Range("B29").Value = 5.237
$ws.Range("B30").Value = 6.208
$ws.Range("B31").Value = 5.578
$ws.Range("E31").Value = 16.267
$ws.Range("A37").Value = -19.92
$ws.Range("E38").Value = 16.489
$ws.Range("B40").Value = 8.962
$ws.Range("E42").Value = 16.586
$ws.Range("E44").Value = 16.93
$ws.Range("E47").Value = 16.271
$ws.Range("B50").Value = 4.946000000000001
$ws.Range("A55").Value = -21.797
$ws.Range("E58").Value = 16.567
$ws.Range("E65").Value = 17.301
$ws.Range("A68").Value = -21.473
$ws.Range("B68").Value = 5.512
$ws.Range("E73").Value = 16.609
$ws.Range("B76").Value = 5.708
$ws.Range("A77").Value = -20.899
$ws.Range("A78").Value = -20.527
$ws.Range("A81").Value = -21.65
$ws.Range("A82").Value = -22.077
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.858000000000001
$ws.Range("E90").Value = 16.363
$ws.Range("E94").Value = 17.828
$ws.Range("E95").Value = 17.564
$ws.Range("B96").Value = 6.692
$ws.Range("B98").Value = 5.355
$ws.Range("B101").Value = 7.527000000000001
$ws.Range("E101").Value = 16.834
$ws.Range("B102").Value = 7.739999999999999
